# The template used Word "field codes" (fldChar begin/instrText/fldChar end)
# to hold the M2Doc directives. The parser was switched to
# TokenIteratorFieldRewriterSplit, which expects the directives to be plain
# literal text delimited by "{" and "}" instead of real Word fields.
#
# This script replaces the two field codes in the second paragraph with
# plain text runs carrying the equivalent "{...}" literal content, while
# leaving the "A paragraph" run between them untouched.

$d = $word.ActiveDocument
$p = $d.Paragraphs(2)

# Range covering the whole paragraph content but excluding the trailing
# paragraph mark, so the mark (and its formatting) is preserved.
$r = $p.Range
[void]$r.MoveEnd(1, -1)

# Remove the existing content (the two fldChar/instrText fields and the
# "A paragraph" run) ...
[void]$r.Delete()

# ... and re-insert it as plain w:t runs, splitting the former field codes
# into individual literal-text runs the same way the original runs were
# split, now wrapped with "{" / "}" instead of begin/end fldChar markers.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{m:</w:t></w:r><w:r><w:t xml:space="preserve">for v </w:t></w:r><w:r><w:t>|</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">self.}</w:t></w:r><w:r><w:t>A paragraph</w:t></w:r><w:r><w:t>{</w:t></w:r><w:r><w:t>m:</w:t></w:r><w:r><w:t xml:space="preserve">endfor}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r.InsertXML($xml)
